# KIBON-112: Zahlungslauf muss pro Gemeinde sein.
# Add a "Gemeinde" column to the ZahlungAuftragPeriode report (Data sheet):
# a new column is inserted before the existing "Institution" detail column
# in the header row (6) and the repeating-row placeholders (7), while the
# unrelated parameter row (4) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 6 (table header): shift existing headers right by one column,
#     starting from the rightmost one so we don't overwrite data we still
#     need to read. (Use Value2 for reads - Value's getter is unreliable
#     here for strings.)
$ws.Range("D6").Value2 = $ws.Range("C6").Value2
$ws.Range("C6").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)

$ws.Range("C6").Value2 = $ws.Range("B6").Value2
$ws.Range("B6").Copy()
$ws.Range("C6").PasteSpecial($xlPasteFormats)

$ws.Range("B6").Value2 = "{gemeindeTitle}"

# --- Row 7 (repeating placeholder row): same right-shift treatment.
$ws.Range("E7").Value2 = $ws.Range("D7").Value2
$ws.Range("D7").Copy()
$ws.Range("E7").PasteSpecial($xlPasteFormats)

$ws.Range("D7").Value2 = $ws.Range("C7").Value2
$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial($xlPasteFormats)

$ws.Range("C7").Value2 = $ws.Range("B7").Value2
$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial($xlPasteFormats)

$ws.Range("B7").Value2 = "{gemeinde}"
$ws.Range("A7").Copy()
$ws.Range("B7").PasteSpecial($xlPasteFormats)

# --- Rows 1-3: widen the title/parameter rows with a (blank, same-styled)
#     cell in column B, matching column A's formatting.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial($xlPasteFormats)

$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial($xlPasteFormats)

$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0
